# Apply the data + selection edits described by the commit diff to "Feuil1".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Row 22 values change: B22 6 -> 2, C22 12 -> 36
$ws.Range("B22").Value = 2
$ws.Range("C22").Value = 36

# Selection moves from A23:C30 (active cell C23) to single cell D25
$ws.Range("D25").Select()
